$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.4240171225616384
$ws.Range("D2").Value = 0.6756706158697217

$ws.Range("C3").Value = 1.297472637530896
$ws.Range("D3").Value = 0.2079096124338125

$ws.Range("C4").Value = 0.3643212262931096
$ws.Range("D4").Value = 0.7190964590771953

$ws.Range("C5").Value = 0.985074431693005
$ws.Range("D5").Value = 0.3352979926210433

$ws.Range("C6").Value = 0.6314105431313786
$ws.Range("D6").Value = 0.5342782959237495

$ws.Range("C7").Value = 0.02712534776445143
$ws.Range("D7").Value = 0.9786043132061861

$ws.Range("C8").Value = 0.2850583328217302
$ws.Range("D8").Value = 0.7782657717969599

$ws.Range("C9").Value = -0.5257199648938949
$ws.Range("D9").Value = 0.6043420937303678

$ws.Range("C10").Value = -0.4024798423297293
$ws.Range("D10").Value = 0.6912149222710582

$ws.Range("C11").Value = 0.2436863133095518
$ws.Range("D11").Value = 0.8097317954739653
